# Edit script for fsar-last-page.docx
# 1) Remove phone number paragraph content, update paragraph-mark rPr
# 2) Rewrite ISSN/ISBN/Cat.No. line, adding report_number_eng bookmark + E-PDF
# 3) Remove the now-duplicate report_number_eng bookmark later in the doc
# Bookmark w:id values are reassigned by Word on save (in document order),
# so we don't need to manage them manually.

$d = $word.ActiveDocument

$xmlPhoneEmail = '<w:p w14:paraId="42F4C367" w14:textId="3C406A34" w:rsidR="00FA1C7B" w:rsidRPr="0033203F" w:rsidRDefault="00FA1C7B" w:rsidP="00CA28AD"><w:pPr><w:pStyle w:val="BodyTextCentered"/><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="0033203F"><w:br/><w:t xml:space="preserve">E-Mail: </w:t></w:r><w:bookmarkStart w:id="3" w:name="email"/><w:r w:rsidR="00D40CF0"><w:t>[</w:t></w:r><w:r w:rsidR="00E834DF"><w:t>email</w:t></w:r><w:r w:rsidR="00D40CF0"><w:t>]</w:t></w:r><w:bookmarkEnd w:id="3"/><w:r w:rsidRPr="0033203F"><w:br/></w:r><w:r w:rsidRPr="0033203F"><w:rPr><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Internet address: </w:t></w:r><w:r w:rsidRPr="0033203F"><w:rPr><w:szCs w:val="22"/></w:rPr><w:fldChar w:fldCharType="begin" w:fldLock="1"/></w:r><w:r w:rsidRPr="0033203F"><w:rPr><w:szCs w:val="22"/></w:rPr><w:instrText>HYPERLINK "http://www.dfo-mpo.gc.ca/csas-sccs/" \o "Fisheries and Oceans Canada / Canadian Science Advisory Secretariat"</w:instrText></w:r><w:r w:rsidRPr="0033203F"><w:rPr><w:szCs w:val="22"/></w:rPr></w:r><w:r w:rsidRPr="0033203F"><w:rPr><w:szCs w:val="22"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidRPr="0033203F"><w:rPr><w:rStyle w:val="Hyperlink"/><w:szCs w:val="22"/></w:rPr><w:t>www.dfo-mpo.gc.ca/csas-sccs/</w:t></w:r></w:p>'
$xmlIssnIsbn = '<w:p w14:paraId="42F9261A" w14:textId="2417AECD" w:rsidR="00EE3838" w:rsidRPr="00D223C8" w:rsidRDefault="00FA1C7B" w:rsidP="00CA28AD"><w:pPr><w:pStyle w:val="BodyTextCentered"/></w:pPr><w:r w:rsidRPr="0033203F"><w:rPr><w:szCs w:val="22"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidR="00EE3838"><w:t xml:space="preserve">ISSN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00EE3838"><w:t>xxxx</w:t></w:r><w:r w:rsidR="00EE3838"><w:t>-</w:t></w:r><w:r w:rsidR="00EE3838"><w:t>xxxx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00EE3838"><w:br/></w:r><w:r w:rsidR="00EE3838"><w:t>ISBN 978-0-660-</w:t></w:r><w:r w:rsidR="00EE3838"><w:rPr><w:color w:val="4F81BD" w:themeColor="accent1"/></w:rPr><w:t>xxxxx-x</w:t></w:r><w:r w:rsidR="00EE3838"><w:tab/><w:t>Cat. No. Fs70-7/2024-</w:t></w:r><w:bookmarkStart w:id="100" w:name="report_number_eng"/><w:r w:rsidR="001D1E53"><w:t>[</w:t></w:r><w:r w:rsidR="00EE3838"><w:t xml:space="preserve">report number </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D1E53"><w:t>eng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D1E53"><w:t>]</w:t></w:r><w:bookmarkEnd w:id="100"/><w:r w:rsidR="00EE3838"><w:t>E-PDF</w:t></w:r><w:r w:rsidR="00EE3838"><w:br/></w:r><w:r w:rsidR="00EE3838" w:rsidRPr="00D223C8"><w:t>© H</w:t></w:r><w:r w:rsidR="00516A3D"><w:t>is</w:t></w:r><w:r w:rsidR="00EE3838" w:rsidRPr="00D223C8"><w:t xml:space="preserve"> Majesty the </w:t></w:r><w:r w:rsidR="00516A3D"><w:t>King</w:t></w:r><w:r w:rsidR="00EE3838" w:rsidRPr="00D223C8"><w:t xml:space="preserve"> in Right of Canada, </w:t></w:r><w:bookmarkStart w:id="4" w:name="copyright_year"/><w:r w:rsidR="00D40CF0"><w:t>[</w:t></w:r><w:r w:rsidR="00E834DF"><w:t>copyright year</w:t></w:r><w:r w:rsidR="00D40CF0"><w:t>]</w:t></w:r><w:bookmarkEnd w:id="4"/></w:p>'
$xmlCitationEng = '<w:p w14:paraId="3E18C7BC" w14:textId="4B1088AC" w:rsidR="00043854" w:rsidRPr="00350C57" w:rsidRDefault="00043854" w:rsidP="00043854"><w:pPr><w:pStyle w:val="citation"/></w:pPr><w:r w:rsidRPr="0092747E"><w:t xml:space="preserve">DFO. </w:t></w:r><w:bookmarkStart w:id="5" w:name="report_year_eng"/><w:r w:rsidR="001D1E53"><w:t>[</w:t></w:r><w:r><w:t>report year</w:t></w:r><w:r w:rsidR="001D1E53"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D1E53"><w:t>eng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D1E53"><w:t>]</w:t></w:r><w:bookmarkEnd w:id="5"/><w:r w:rsidRPr="0092747E"><w:t xml:space="preserve">.  </w:t></w:r><w:bookmarkStart w:id="6" w:name="report_title_eng"/><w:r w:rsidR="001D1E53"><w:t>[</w:t></w:r><w:r><w:t>report title</w:t></w:r><w:r w:rsidR="001D1E53"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D1E53"><w:t>eng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D1E53"><w:t>]</w:t></w:r><w:bookmarkEnd w:id="6"/><w:r w:rsidRPr="0092747E"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="008E1A63" w:rsidRPr="0092747E"><w:t xml:space="preserve">DFO Can. Sci. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="008E1A63" w:rsidRPr="00350C57"><w:t>Advis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="008E1A63" w:rsidRPr="00350C57"><w:t xml:space="preserve">. Sec. Sci. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="008E1A63" w:rsidRPr="00350C57"><w:t>Advis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="008E1A63" w:rsidRPr="00350C57"><w:t xml:space="preserve">. Rep. </w:t></w:r><w:bookmarkStart w:id="7" w:name="report_year_eng2"/><w:r w:rsidR="001D1E53"><w:t>[</w:t></w:r><w:r><w:t>report year</w:t></w:r><w:r w:rsidR="001D1E53"><w:t xml:space="preserve"> eng</w:t></w:r><w:r w:rsidR="00776A19"><w:t>2</w:t></w:r><w:r w:rsidR="001D1E53"><w:t>]</w:t></w:r><w:bookmarkEnd w:id="7"/><w:r><w:t>/</w:t></w:r><w:r w:rsidR="001D1E53"><w:t>[</w:t></w:r><w:r><w:t>report number</w:t></w:r><w:r w:rsidR="001D1E53"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D1E53"><w:t>eng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D1E53"><w:t>]</w:t></w:r><w:r w:rsidRPr="00350C57"><w:t>.</w:t></w:r></w:p>'

$d.Paragraphs(3).Range.InsertXML($xmlPhoneEmail)
$d.Paragraphs(4).Range.InsertXML($xmlIssnIsbn)
$d.Paragraphs(7).Range.InsertXML($xmlCitationEng)

Write-Output "done"
